# Insert a new "INDEX (DO NOT MODIFY)" column at the very left of the
# sheet (shifting every existing column one place to the right, along
# with their data/styles/widths and the dataValidation range that
# referenced the old last column), then populate the new column and
# uppercase the header row (except the very last header, which keeps
# its original casing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRows = $ws.UsedRange.Rows.Count

# 1) Insert a brand-new column before column A. Excel shifts every
#    existing column (data + styles) one place to the right and
#    updates dependent ranges (e.g. the dataValidation sqref) for us.
$ws.Columns.Item(1).Insert()

# 2) The freshly inserted column has no formatting of its own yet;
#    copy it from the column immediately to its right (which holds
#    what used to be column A) so the new column matches it exactly
#    (same header style, same per-row data style).
$ws.Range("B1:B$usedRows").Copy()
$ws.Range("A1:A$usedRows").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Give the new column the same width (23 characters) used in the
# target file. ColumnWidth is expressed in characters and Excel adds
# ~5/6 of a character of internal padding when storing it, so back
# that padding out to land exactly on a stored width of 23.
$ws.Columns.Item(1).ColumnWidth = 23 - 5/6

# 3) Header + data for the new column.
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$indexValues = @(9, 10, 11, 134, 135, 136, 180, 181, 182, 183, 184, 212, 213, 214, 215, 216, 217)
for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $indexValues[$i]
}

# 4) Uppercase the rest of the header row, except the very last column
#    (the "Status as of ..." column), which is left unchanged.
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column
for ($c = 2; $c -lt $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $txt = $cell.Text
    $cell.Value = $txt.ToUpper()
}
